{"js": "// Applies the LOM3073.docx content revisions (2023-04-12 build) by\n// locating each original run's exact text and replacing it in place.\nconst body = context.document.body;\n\n// Change 1\n{\n  const results = body.search(\"Ativa\u00e7\u00e3o: 01/01/2022\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(`Change 1: expected 1 match, found ${results.items.length}`);\n  }\n  results.items[0].insertText(\"Ativa\u00e7\u00e3o: 01/01/2023\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Change 2\n{\n  const results = body.search(\"1 - Capacitar o aluno a entender e discutir as teorias e fen\u00f4menos envolvidos no processamento cer\u00e2mico, 2 - Conhecer as mat\u00e9rias-primas naturais e sint\u00e9ticas, bem como suas ocorr\u00eancias e propriedade, 3 - Capacitar o aluno a utilizar instrumentos e m\u00e9todos que visam transformar mat\u00e9rias-primas cer\u00e2micas em produtos\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(`Change 2: expected 1 match, found ${results.items.length}`);\n  }\n  results.items[0].insertText(\"Esta disciplina faz parte da forma\u00e7\u00e3o do engenheiro de materiais, inserida na grande \u00e1rea \u201cCER\u00c2MICA\u201d, contribuindo para gerar compet\u00eancias gerais e espec\u00edficas.Capacitar os alunos a desenvolverem uma vis\u00e3o integrada da viabilidade t\u00e9cnico-econ\u00f4mica-ambiental das principais tecnologias envolvidas no processamento de cer\u00e2micas.Incentivar trabalhos em grupo, com \u00eanfase na vis\u00e3o integrada sobre os aspectos abordados na disciplina.Promover a comunica\u00e7\u00e3o nas formas escrita, oral e gr\u00e1fica, al\u00e9m de trabalhos em grupos.Relacionar esta disciplina com outras da grade do curso, tanto com as de forma\u00e7\u00e3o espec\u00edfica quanto \u00e0s de forma\u00e7\u00e3o geral.Apresentar os principais equipamentos envolvidos no processamento de cer\u00e2micas preparando os alunos para a disciplina \u201cProcessamento de Cer\u00e2micas Experimental\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Change 3\n{\n  const results = body.search(\"1 \u2013 Introdu\u00e7\u00e3o, 2 - Mat\u00e9rias primas naturais e sint\u00e9ticas, 3 - Preparo de massas cer\u00e2micas, 4 \u2013 Conforma\u00e7\u00e3o, 5 - Vari\u00e1veis cr\u00edticas no controle do processamento, 6 \u2013 Testes experimentais\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(`Change 3: expected 1 match, found ${results.items.length}`);\n  }\n  results.items[0].insertText(\"1 \u2013 Introdu\u00e7\u00e3o, 2 \u2013 Principais mat\u00e9rias-primas naturais e sint\u00e9ticas, 3 - Preparo de massas cer\u00e2micas, 4 \u2013 Conforma\u00e7\u00e3o, 5 \u2013 Queima (sintereiza\u00e7\u00e3o), 6 - Vari\u00e1veis cr\u00edticas no controle do processamento, 7 \u2013 Usinagem de materiais cer\u00e2micos. 8 \u2013 Projeto\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Change 4\n{\n  const results = body.search(\"1 \u2013 Introdu\u00e7\u00e3o - conceito e defini\u00e7\u00f5es. Hist\u00f3ria; 2 \u2013 Mat\u00e9rias primas - naturais pl\u00e1sticas e n\u00e3o pl\u00e1sticas: origem, ocorr\u00eancia, propriedades f\u00edsico-qu\u00edmicas, beneficiamento; 3 \u2013 Mat\u00e9rias-primas sint\u00e9ticas: s\u00ednteses e processos qu\u00edmicos; 4 \u2013 Caracteriza\u00e7\u00e3o f\u00edsico-qu\u00edmica dos p\u00f3s cer\u00e2micos; 5 \u2013 Aditivos: plastificantes, defloculantes, lubrificantes e agentes suspensores; 6 \u2013 Preparo de massas cer\u00e2micas - tipo de massas: branca, vermelha, refrat\u00e1ria e especiais; preparo e propriedades reol\u00f3gicas de pastas e suspens\u00f5es cer\u00e2micas (moagem, tipos de moinhos, carregamento e contamina\u00e7\u00e3o - m\u00e9todos de aglomera\u00e7\u00e3o: \u201cspray drier\u201de disco (misturador Eirich); 7 \u2013 Conforma\u00e7\u00e3o via seca: teoria, tipos; equipamentos, etapas do ciclo de prensagem e defeitos nas pe\u00e7as; 8 \u2013 Conforma\u00e7\u00e3o via l\u00edquida e pastas: teoria, estabilidade e reologia das suspens\u00f5es, diagrama de Atterberg, moldes, equipamentos, etapas do ciclo e defeitos caracter\u00edsticos da forma\u00e7\u00e3o; 9 - Vari\u00e1veis cr\u00edticas no controle do processamento - influ\u00eancia das caracter\u00edsticas dos materiais e dos equipamentos no estabelecimento das condi\u00e7\u00f5es de conforma\u00e7\u00e3o. 10 \u2013 Fundamentos e processos inovadores de conforma\u00e7\u00e3o de materiais cer\u00e2micos. 11 - Testes experimentais.\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(`Change 4: expected 1 match, found ${results.items.length}`);\n  }\n  results.items[0].insertText(\"1 \u2013 Introdu\u00e7\u00e3o: defini\u00e7\u00f5es, setores cer\u00e2micos, aplica\u00e7\u00f5es e fontes de divulga\u00e7\u00e3o da \u00e1rea de cer\u00e2micas. 2 - Principais mat\u00e9rias-primas naturais e sint\u00e9ticas: ocorrencias naturais e beneficiamento e s\u00ednteses de cer\u00e2micas avan\u00e7adas. 3 - Preparo de massas cer\u00e2micas: formula\u00e7\u00e3o de composi\u00e7\u00f5es cer\u00e2micas com e sem utiliza\u00e7\u00e3o de diagramas de fases, reologia das barbotinas e pastas. 4 \u2013 Conforma\u00e7\u00e3o:  equipamentos utilizados na conforma\u00e7\u00e3o de cer\u00e2micas tradicionais e t\u00e9cnicas, defeitos e problemas na conforma\u00e7\u00e3o, m\u00e9todos de conforma\u00e7\u00e3o (colagem de barbotina, prensagem, extrus\u00e3o, inje\u00e7\u00e3o). 5 \u2013 Queima (sintereiza\u00e7\u00e3o): curva de queima, eventos pr\u00e9-sinteriza\u00e7\u00e3o, sinteriza\u00e7\u00e3o, mecanismos de sinteriza\u00e7\u00e3o, equipamentos, sinteriza\u00e7\u00e3o r\u00e1pida, microestrutura (controle microestrutural, rela\u00e7\u00e3o microestrutura x propriedades), 6 - Vari\u00e1veis cr\u00edticas no controle do processamento: avaliadas em cada etapa do processamento. 7 - Usinagem de materiais cer\u00e2micos: usinagem a verde e ap\u00f3s sinteriza\u00e7\u00e3o, defeitos superficiais introduzidos, acabamento.8 \u2013 Projeto: Desenvolvimento de produtos cer\u00e2micos levando em conta aspectos de inova\u00e7\u00e3o, sustentabilidade, social e ec\u00f4n\u00f4mico. Este t\u00f3pico dever\u00e1 ser desenvolvido em grupo.\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Change 5\n{\n  const results = body.search(\"a) Duas provas escritas (P1 e P2, com peso 1)b) Relat\u00f3rios sobre os testes experimentais: soma das notas dos relat\u00f3rios divido pelo n\u00famero de relat\u00f3rios (SR), com peso 1.\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(`Change 5: expected 1 match, found ${results.items.length}`);\n  }\n  results.items[0].insertText(\"O aprendizado te\u00f3rico ser\u00e1 avaliado por meio de provas escritas presenciais e os conhecimentos gerais e espec\u00edficos ser\u00e3o avaliados pela monografia e apresenta\u00e7\u00e3o do projeto, e participa\u00e7\u00e3o nas discuss\u00f5es individuais e em grupos.\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Change 6\n{\n  const results = body.search(\"Ser\u00e3o aplicadas duas avalia\u00e7\u00f5es escritas (P1 e P2, com peso 1). A nota final ser\u00e3o calculada pela equa\u00e7\u00e3oNF = (P1+P2+MR)/3. NF igual ou superior a 5: aprova\u00e7\u00e3o direta. NF entre 3 e 4,9: recupera\u00e7\u00e3o. NF inferior a 3: reprova\u00e7\u00e3o direta.\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(`Change 6: expected 1 match, found ${results.items.length}`);\n  }\n  results.items[0].insertText(\"(a) duas provas escritas (P1 e P2) valendo de zero a dez, (b) nota de projeto (NP): monografia (variando de zero a seis), (c) apresenta\u00e7\u00e3o oral, (AO), (variando de zero a quatro), (d) participa\u00e7\u00e3o nas aulas (PA), valendo de zero a dez).A nota final (NF) ser\u00e1 calculada pela equa\u00e7\u00e3o:NF = 0,9 x [(P1+P2+NP)/3]+0,1 (PA) NF igual ou superior a 5: aprova\u00e7\u00e3o direta. NF entre 3 e 4,9: recupera\u00e7\u00e3o. NF inferior a 3: reprova\u00e7\u00e3o direta.\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Change 7\n{\n  const results = body.search(\"1. Santos, P.S. Tecnologia de Argilas, vol. 2, EDUSP, 1975 e 1989,2. Dispers\u00e3o e empacotamento de part\u00edculas, Fazendo Arte Editorial. Ivone R de Oliveira e co-autores, 2000,3. Norton, F.H. Introdu\u00e7\u00e3o \u00e0 Tecnologia Cer\u00e2mica, Ed. Edgard Blucher, 1973,4. Kingery, W.D. Introduction to Ceramics, John Wiley, 1970 e 1976, 2nd Edition, 5. Reed, J.S. Principles of Ceramics Processing, John Wiley, 1988, 6. Rahaman, M. N. Ceramic Processing and Sintering. 1st Edition, 1993, 7.Van Vlack, L.M. Propriedades dos Materiais Cer\u00e2micos, Ed. Edgard Blucher, 1973, 8. Ceramic Materials: Science and Engineering, C. Barry Carter, M. Grant Norton  2nd ed., 2013,9. Fundamentals of Ceramic Powder Processing and Synthesis: Terry A. Ring10. Artigos da literatura especializada\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(`Change 7: expected 1 match, found ${results.items.length}`);\n  }\n  results.items[0].insertText(\"1. Santos, P.S. Tecnologia de Argilas, vol. 2, EDUSP, 1975 e 1989, 2. Dispers\u00e3o e empacotamento de part\u00edculas, Fazendo Arte Editorial. Ivone R de Oliveira e co-autores, 2000, 3. Norton, F.H. Introdu\u00e7\u00e3o \u00e0 Tecnologia Cer\u00e2mica, Ed. Edgard Blucher, 1973, 4. Kingery, W.D. Introduction to Ceramics, John Wiley, 1970 e 1976, 2nd Edition, 5. Reed, J.S. Principles of Ceramics Processing, John Wiley, 1988, 6. Rahaman, M. N. Ceramic Processing and Sintering. 1st Edition, 1993, 7.Van Vlack, L.M. Propriedades dos Materiais Cer\u00e2micos, Ed. Edgard Blucher, 1973, 8. Ceramic Materials: Science and Engineering, C. Barry Carter, M. Grant Norton 2nd ed., 2013, 9. Fundamentals of Ceramic Powder Processing and Synthesis: Terry A. Ring, 11. Setz, L.F. G. O Processamento Cer\u00e2mico sem Mist\u00e9rio. 1\u00aa edi\u00e7\u00e3o, 256 p\u00e1ginas, Edgard Bl\u00fccher, 2019, 12. M. F. Ashby, D.R. H. Jones, Engenharia de Materiais, Volume II, 3\u00b0 edi\u00e7\u00e3o, Elsevier, p.436, 2007, 13. Artigos da literatura especializada\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Applies the LOM3073.docx content revisions (2023-04-12 build) by\n# locating each original run's exact text and replacing it in place.\n$d = $word.ActiveDocument\n\n# Change 1\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = 'Ativa\u00e7\u00e3o: 01/01/2022'\n$find.Replacement.Text = 'Ativa\u00e7\u00e3o: 01/01/2023'\n$ok = $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\nif (-not $ok) { throw \"Change 1: Find.Execute did not find/replace the expected text\" }\n\n# Change 2\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '1 - Capacitar o aluno a entender e discutir as teorias e fen\u00f4menos envolvidos no processamento cer\u00e2mico, 2 - Conhecer as mat\u00e9rias-primas naturais e sint\u00e9ticas, bem como suas ocorr\u00eancias e propriedade, 3 - Capacitar o aluno a utilizar instrumentos e m\u00e9todos que visam transformar mat\u00e9rias-primas cer\u00e2micas em produtos'\n$find.Replacement.Text = 'Esta disciplina faz parte da forma\u00e7\u00e3o do engenheiro de materiais, inserida na grande \u00e1rea \u201cCER\u00c2MICA\u201d, contribuindo para gerar compet\u00eancias gerais e espec\u00edficas.Capacitar os alunos a desenvolverem uma vis\u00e3o integrada da viabilidade t\u00e9cnico-econ\u00f4mica-ambiental das principais tecnologias envolvidas no processamento de cer\u00e2micas.Incentivar trabalhos em grupo, com \u00eanfase na vis\u00e3o integrada sobre os aspectos abordados na disciplina.Promover a comunica\u00e7\u00e3o nas formas escrita, oral e gr\u00e1fica, al\u00e9m de trabalhos em grupos.Relacionar esta disciplina com outras da grade do curso, tanto com as de forma\u00e7\u00e3o espec\u00edfica quanto \u00e0s de forma\u00e7\u00e3o geral.Apresentar os principais equipamentos envolvidos no processamento de cer\u00e2micas preparando os alunos para a disciplina \u201cProcessamento de Cer\u00e2micas Experimental'\n$ok = $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\nif (-not $ok) { throw \"Change 2: Find.Execute did not find/replace the expected text\" }\n\n# Change 3\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '1 \u2013 Introdu\u00e7\u00e3o, 2 - Mat\u00e9rias primas naturais e sint\u00e9ticas, 3 - Preparo de massas cer\u00e2micas, 4 \u2013 Conforma\u00e7\u00e3o, 5 - Vari\u00e1veis cr\u00edticas no controle do processamento, 6 \u2013 Testes experimentais'\n$find.Replacement.Text = '1 \u2013 Introdu\u00e7\u00e3o, 2 \u2013 Principais mat\u00e9rias-primas naturais e sint\u00e9ticas, 3 - Preparo de massas cer\u00e2micas, 4 \u2013 Conforma\u00e7\u00e3o, 5 \u2013 Queima (sintereiza\u00e7\u00e3o), 6 - Vari\u00e1veis cr\u00edticas no controle do processamento, 7 \u2013 Usinagem de materiais cer\u00e2micos. 8 \u2013 Projeto'\n$ok = $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\nif (-not $ok) { throw \"Change 3: Find.Execute did not find/replace the expected text\" }\n\n# Change 4\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '1 \u2013 Introdu\u00e7\u00e3o - conceito e defini\u00e7\u00f5es. Hist\u00f3ria; 2 \u2013 Mat\u00e9rias primas - naturais pl\u00e1sticas e n\u00e3o pl\u00e1sticas: origem, ocorr\u00eancia, propriedades f\u00edsico-qu\u00edmicas, beneficiamento; 3 \u2013 Mat\u00e9rias-primas sint\u00e9ticas: s\u00ednteses e processos qu\u00edmicos; 4 \u2013 Caracteriza\u00e7\u00e3o f\u00edsico-qu\u00edmica dos p\u00f3s cer\u00e2micos; 5 \u2013 Aditivos: plastificantes, defloculantes, lubrificantes e agentes suspensores; 6 \u2013 Preparo de massas cer\u00e2micas - tipo de massas: branca, vermelha, refrat\u00e1ria e especiais; preparo e propriedades reol\u00f3gicas de pastas e suspens\u00f5es cer\u00e2micas (moagem, tipos de moinhos, carregamento e contamina\u00e7\u00e3o - m\u00e9todos de aglomera\u00e7\u00e3o: \u201cspray drier\u201de disco (misturador Eirich); 7 \u2013 Conforma\u00e7\u00e3o via seca: teoria, tipos; equipamentos, etapas do ciclo de prensagem e defeitos nas pe\u00e7as; 8 \u2013 Conforma\u00e7\u00e3o via l\u00edquida e pastas: teoria, estabilidade e reologia das suspens\u00f5es, diagrama de Atterberg, moldes, equipamentos, etapas do ciclo e defeitos caracter\u00edsticos da forma\u00e7\u00e3o; 9 - Vari\u00e1veis cr\u00edticas no controle do processamento - influ\u00eancia das caracter\u00edsticas dos materiais e dos equipamentos no estabelecimento das condi\u00e7\u00f5es de conforma\u00e7\u00e3o. 10 \u2013 Fundamentos e processos inovadores de conforma\u00e7\u00e3o de materiais cer\u00e2micos. 11 - Testes experimentais.'\n$find.Replacement.Text = '1 \u2013 Introdu\u00e7\u00e3o: defini\u00e7\u00f5es, setores cer\u00e2micos, aplica\u00e7\u00f5es e fontes de divulga\u00e7\u00e3o da \u00e1rea de cer\u00e2micas. 2 - Principais mat\u00e9rias-primas naturais e sint\u00e9ticas: ocorrencias naturais e beneficiamento e s\u00ednteses de cer\u00e2micas avan\u00e7adas. 3 - Preparo de massas cer\u00e2micas: formula\u00e7\u00e3o de composi\u00e7\u00f5es cer\u00e2micas com e sem utiliza\u00e7\u00e3o de diagramas de fases, reologia das barbotinas e pastas. 4 \u2013 Conforma\u00e7\u00e3o:  equipamentos utilizados na conforma\u00e7\u00e3o de cer\u00e2micas tradicionais e t\u00e9cnicas, defeitos e problemas na conforma\u00e7\u00e3o, m\u00e9todos de conforma\u00e7\u00e3o (colagem de barbotina, prensagem, extrus\u00e3o, inje\u00e7\u00e3o). 5 \u2013 Queima (sintereiza\u00e7\u00e3o): curva de queima, eventos pr\u00e9-sinteriza\u00e7\u00e3o, sinteriza\u00e7\u00e3o, mecanismos de sinteriza\u00e7\u00e3o, equipamentos, sinteriza\u00e7\u00e3o r\u00e1pida, microestrutura (controle microestrutural, rela\u00e7\u00e3o microestrutura x propriedades), 6 - Vari\u00e1veis cr\u00edticas no controle do processamento: avaliadas em cada etapa do processamento. 7 - Usinagem de materiais cer\u00e2micos: usinagem a verde e ap\u00f3s sinteriza\u00e7\u00e3o, defeitos superficiais introduzidos, acabamento.8 \u2013 Projeto: Desenvolvimento de produtos cer\u00e2micos levando em conta aspectos de inova\u00e7\u00e3o, sustentabilidade, social e ec\u00f4n\u00f4mico. Este t\u00f3pico dever\u00e1 ser desenvolvido em grupo.'\n$ok = $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\nif (-not $ok) { throw \"Change 4: Find.Execute did not find/replace the expected text\" }\n\n# Change 5\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = 'a) Duas provas escritas (P1 e P2, com peso 1)b) Relat\u00f3rios sobre os testes experimentais: soma das notas dos relat\u00f3rios divido pelo n\u00famero de relat\u00f3rios (SR), com peso 1.'\n$find.Replacement.Text = 'O aprendizado te\u00f3rico ser\u00e1 avaliado por meio de provas escritas presenciais e os conhecimentos gerais e espec\u00edficos ser\u00e3o avaliados pela monografia e apresenta\u00e7\u00e3o do projeto, e participa\u00e7\u00e3o nas discuss\u00f5es individuais e em grupos.'\n$ok = $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\nif (-not $ok) { throw \"Change 5: Find.Execute did not find/replace the expected text\" }\n\n# Change 6\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = 'Ser\u00e3o aplicadas duas avalia\u00e7\u00f5es escritas (P1 e P2, com peso 1). A nota final ser\u00e3o calculada pela equa\u00e7\u00e3oNF = (P1+P2+MR)/3. NF igual ou superior a 5: aprova\u00e7\u00e3o direta. NF entre 3 e 4,9: recupera\u00e7\u00e3o. NF inferior a 3: reprova\u00e7\u00e3o direta.'\n$find.Replacement.Text = '(a) duas provas escritas (P1 e P2) valendo de zero a dez, (b) nota de projeto (NP): monografia (variando de zero a seis), (c) apresenta\u00e7\u00e3o oral, (AO), (variando de zero a quatro), (d) participa\u00e7\u00e3o nas aulas (PA), valendo de zero a dez).A nota final (NF) ser\u00e1 calculada pela equa\u00e7\u00e3o:NF = 0,9 x [(P1+P2+NP)/3]+0,1 (PA) NF igual ou superior a 5: aprova\u00e7\u00e3o direta. NF entre 3 e 4,9: recupera\u00e7\u00e3o. NF inferior a 3: reprova\u00e7\u00e3o direta.'\n$ok = $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\nif (-not $ok) { throw \"Change 6: Find.Execute did not find/replace the expected text\" }\n\n# Change 7\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = '1. Santos, P.S. Tecnologia de Argilas, vol. 2, EDUSP, 1975 e 1989,2. Dispers\u00e3o e empacotamento de part\u00edculas, Fazendo Arte Editorial. Ivone R de Oliveira e co-autores, 2000,3. Norton, F.H. Introdu\u00e7\u00e3o \u00e0 Tecnologia Cer\u00e2mica, Ed. Edgard Blucher, 1973,4. Kingery, W.D. Introduction to Ceramics, John Wiley, 1970 e 1976, 2nd Edition, 5. Reed, J.S. Principles of Ceramics Processing, John Wiley, 1988, 6. Rahaman, M. N. Ceramic Processing and Sintering. 1st Edition, 1993, 7.Van Vlack, L.M. Propriedades dos Materiais Cer\u00e2micos, Ed. Edgard Blucher, 1973, 8. Ceramic Materials: Science and Engineering, C. Barry Carter, M. Grant Norton  2nd ed., 2013,9. Fundamentals of Ceramic Powder Processing and Synthesis: Terry A. Ring10. Artigos da literatura especializada'\n$find.Replacement.Text = '1. Santos, P.S. Tecnologia de Argilas, vol. 2, EDUSP, 1975 e 1989, 2. Dispers\u00e3o e empacotamento de part\u00edculas, Fazendo Arte Editorial. Ivone R de Oliveira e co-autores, 2000, 3. Norton, F.H. Introdu\u00e7\u00e3o \u00e0 Tecnologia Cer\u00e2mica, Ed. Edgard Blucher, 1973, 4. Kingery, W.D. Introduction to Ceramics, John Wiley, 1970 e 1976, 2nd Edition, 5. Reed, J.S. Principles of Ceramics Processing, John Wiley, 1988, 6. Rahaman, M. N. Ceramic Processing and Sintering. 1st Edition, 1993, 7.Van Vlack, L.M. Propriedades dos Materiais Cer\u00e2micos, Ed. Edgard Blucher, 1973, 8. Ceramic Materials: Science and Engineering, C. Barry Carter, M. Grant Norton 2nd ed., 2013, 9. Fundamentals of Ceramic Powder Processing and Synthesis: Terry A. Ring, 11. Setz, L.F. G. O Processamento Cer\u00e2mico sem Mist\u00e9rio. 1\u00aa edi\u00e7\u00e3o, 256 p\u00e1ginas, Edgard Bl\u00fccher, 2019, 12. M. F. Ashby, D.R. H. Jones, Engenharia de Materiais, Volume II, 3\u00b0 edi\u00e7\u00e3o, Elsevier, p.436, 2007, 13. Artigos da literatura especializada'\n$ok = $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\nif (-not $ok) { throw \"Change 7: Find.Execute did not find/replace the expected text\" }\n"}
